$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Wednesday (D13) hours and fill in Thursday (E13) / Friday (F13) hours
# for the week-13 row, matching the "fixed buy/sell bug" timesheet update.
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 1

# Move the active selection to G17, as recorded in the saved view state.
$ws.Range("G17").Select()

$wb.Save()
